$wb = $excel.ActiveWorkbook

# "Datos a Cargar": the sheet shipped with three rows of sample/demo data
# (rows 2-4) under the header row. Remove the sample values (formatting /
# styles stay in place) so the sheet is ready for real data entry, then
# move the active selection down to A8.
$ws1 = $wb.Worksheets.Item("Datos a Cargar")
$ws1.Range("A2:P4").ClearContents() | Out-Null
$ws1.Range("A8").Select() | Out-Null

# "Guia": scroll back to the top-left corner of the sheet and select A2
# (previously the view was scrolled over to column I with P4 selected).
$ws2 = $wb.Worksheets.Item("Guia")
$ws2.Select() | Out-Null
$ws2.Range("A2").Select() | Out-Null

# Restore "Datos a Cargar" as the active/visible tab and selection, matching
# the original workbook state (it was the tab shown when the file is opened).
$ws1.Select() | Out-Null
$ws1.Range("A8").Select() | Out-Null

# Tab-bar split ratio (cosmetic window setting, between the sheet tabs and
# the horizontal scroll bar).
$excel.ActiveWindow.TabRatio = 0.504
